$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 6414050.5
$ws.Range("I9").Value = 7576417
$ws.Range("K9").Value = 7576417
$ws.Range("M9").Value = -7576248
$ws.Range("H32").Value = 19231932
$ws.Range("I32").Value = 1058.4286
$ws.Range("J32").Value = 26316992
$ws.Range("K32").Value = 1058.4286
$ws.Range("L32").Value = 26316992
$ws.Range("M32").Value = -732.4286
$ws.Range("N32").Value = -26317644
$ws.Range("H98").Value = 2334.4102
$ws.Range("I98").Value = 2008.9375
$ws.Range("J98").Value = 3822.2856
$ws.Range("K98").Value = 2008.9375
$ws.Range("L98").Value = 3822.2856
$ws.Range("M98").Value = -510.9375
$ws.Range("N98").Value = -6818.2856
$ws.Range("H122").Value = 2334.4102
$ws.Range("I122").Value = 2008.9375
$ws.Range("J122").Value = 3822.2856
$ws.Range("K122").Value = 6026.8125
$ws.Range("L122").Value = 11466.8568
$ws.Range("M122").Value = -3576.8125
$ws.Range("N122").Value = -16366.8568
$ws.Range("H132").Value = 1589.7778
$ws.Range("I132").Value = 1308.3572
$ws.Range("J132").Value = 2574.75
$ws.Range("K132").Value = 3925.0716
$ws.Range("L132").Value = 7724.25
$ws.Range("M132").Value = -1395.0716
$ws.Range("N132").Value = -12784.25
$ws.Range("H137").Value = 1376535.8
$ws.Range("I137").Value = 6581.8125
$ws.Range("J137").Value = 2445768.2
$ws.Range("K137").Value = 19745.4375
$ws.Range("L137").Value = 7337304.600000001
$ws.Range("M137").Value = -17195.4375
$ws.Range("N137").Value = -7342404.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 88458.836
$ws.Range("I45").Value = 96000.55
$ws.Range("J45").Value = 5500
$ws.Range("K45").Value = 96000.55
$ws.Range("L45").Value = 5500
$ws.Range("M45").Value = -95623.55
$ws.Range("N45").Value = -6254
$ws.Range("H122").Value = 1869.25
$ws.Range("I122").Value = 1410.3334
$ws.Range("K122").Value = 4231.0002
$ws.Range("M122").Value = -1781.0002
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null
$ws.Range("H132").Value = 2577.1562
$ws.Range("I132").Value = 1653.069
$ws.Range("J132").Value = 11510
$ws.Range("K132").Value = 4959.207
$ws.Range("L132").Value = 34530
$ws.Range("M132").Value = -2429.207
$ws.Range("N132").Value = -39590

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3433.6843
$ws.Range("I86").Value = 2595.739
$ws.Range("J86").Value = 4718.533
$ws.Range("K86").Value = 2595.739
$ws.Range("L86").Value = 4718.533
$ws.Range("M86").Value = -1472.739
$ws.Range("N86").Value = -6964.533
$ws.Range("H89").Value = 3433.6843
$ws.Range("I89").Value = 2595.739
$ws.Range("J89").Value = 4718.533
$ws.Range("K89").Value = 12978.695
$ws.Range("L89").Value = 23592.665
$ws.Range("M89").Value = -7362.695
$ws.Range("N89").Value = -34824.665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5134.131
$ws.Range("I31").Value = 4116.3335
$ws.Range("J31").Value = 5303.7637
$ws.Range("K31").Value = 4116.3335
$ws.Range("L31").Value = 5303.7637
$ws.Range("M31").Value = -3821.3335
$ws.Range("N31").Value = -5893.7637
$ws.Range("H34").Value = 5134.131
$ws.Range("I34").Value = 4116.3335
$ws.Range("J34").Value = 5303.7637
$ws.Range("K34").Value = 4116.3335
$ws.Range("L34").Value = 5303.7637
$ws.Range("M34").Value = -3914.3335
$ws.Range("N34").Value = -5707.7637
$ws.Range("H132").Value = 4026.9333
$ws.Range("I132").Value = 4184.923
$ws.Range("K132").Value = 12554.769
$ws.Range("M132").Value = -10024.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 9522.684999999999
$ws.Range("I3").Value = 7577.875
$ws.Range("J3").Value = 19895
$ws.Range("K3").Value = 22733.625
$ws.Range("L3").Value = 59685
$ws.Range("M3").Value = -22621.625
$ws.Range("N3").Value = -59909
$ws.Range("H14").Value = 3455.1428
$ws.Range("I14").Value = 3455.1428
$ws.Range("K14").Value = 10365.4284
$ws.Range("M14").Value = -10192.4284
$ws.Range("H68").Value = 2854.0952
$ws.Range("I68").Value = 2221.3572
$ws.Range("K68").Value = 6664.071599999999
$ws.Range("M68").Value = -5853.071599999999
$ws.Range("H71").Value = 2854.0952
$ws.Range("I71").Value = 2221.3572
$ws.Range("K71").Value = 19992.2148
$ws.Range("M71").Value = -15936.2148
$ws.Range("H107").Value = 100000980
$ws.Range("I107").Value = 535.5714
$ws.Range("J107").Value = 333335330
$ws.Range("K107").Value = 1606.7142
$ws.Range("L107").Value = 1000005990
$ws.Range("M107").Value = 313.2857999999999
$ws.Range("N107").Value = -1000009830
$ws.Range("H141").Value = 3007.25
$ws.Range("I141").Value = 3007.25
$ws.Range("K141").Value = 9021.75
$ws.Range("M141").Value = -3841.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 7944
$ws.Range("I5").Value = 7944
$ws.Range("K5").Value = 7944
$ws.Range("M5").Value = -7832
$ws.Range("H70").Value = 5192.8335
$ws.Range("I70").Value = 5211.273
$ws.Range("J70").Value = 4990
$ws.Range("K70").Value = 5211.273
$ws.Range("L70").Value = 4990
$ws.Range("M70").Value = -4941.273
$ws.Range("N70").Value = -5530
$ws.Range("H73").Value = 5192.8335
$ws.Range("I73").Value = 5211.273
$ws.Range("J73").Value = 4990
$ws.Range("K73").Value = 5211.273
$ws.Range("L73").Value = 4990
$ws.Range("M73").Value = -4275.273
$ws.Range("N73").Value = -6862
$ws.Range("H113").Value = 4177.0454
$ws.Range("I113").Value = 4029.4443
$ws.Range("J113").Value = 4841.25
$ws.Range("K113").Value = 4029.4443
$ws.Range("L113").Value = 4841.25
$ws.Range("M113").Value = -1859.4443
$ws.Range("N113").Value = -9181.25
$ws.Range("H132").Value = 918326.9
$ws.Range("I132").Value = 2793.8635
$ws.Range("J132").Value = 4275281.5
$ws.Range("K132").Value = 8381.5905
$ws.Range("L132").Value = 12825844.5
$ws.Range("M132").Value = -5851.5905
$ws.Range("N132").Value = -12830904.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8600
$ws.Range("I16").Value = 8450
$ws.Range("J16").Value = 8750
$ws.Range("K16").Value = 8450
$ws.Range("L16").Value = 8750
$ws.Range("M16").Value = -8280
$ws.Range("N16").Value = -9090
$ws.Range("H99").Value = 49048.1
$ws.Range("I99").Value = 32275.666
$ws.Range("K99").Value = 32275.666
$ws.Range("M99").Value = -29280.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").Value = $null
$ws.Range("H81").Value = 204299.6
$ws.Range("I81").Value = 5499.6665
$ws.Range("J81").Value = 502499.5
$ws.Range("K81").Value = 10999.333
$ws.Range("L81").Value = 1004999
$ws.Range("M81").Value = -9938.333000000001
$ws.Range("N81").Value = -1007121
$ws.Range("H84").Value = 204299.6
$ws.Range("I84").Value = 5499.6665
$ws.Range("J84").Value = 502499.5
$ws.Range("K84").Value = 54996.665
$ws.Range("L84").Value = 5024995
$ws.Range("M84").Value = -49692.665
$ws.Range("N84").Value = -5035603
$ws.Range("H113").Value = 780.0714
$ws.Range("I113").Value = 358.4
$ws.Range("J113").Value = 1834.25
$ws.Range("K113").Value = 1075.2
$ws.Range("L113").Value = 5502.75
$ws.Range("M113").Value = 1094.8
$ws.Range("N113").Value = -9842.75
$ws.Range("H132").Value = 1517.2693
$ws.Range("I132").Value = 1356.4584
$ws.Range("K132").Value = 4069.3752
$ws.Range("M132").Value = -1539.3752
$ws.Range("H141").Value = 191989.4
$ws.Range("J141").Value = 191989.4
$ws.Range("L141").Value = 191989.4
$ws.Range("N141").Value = -202349.4
